# LCV duzeltmeleri yapildi, HDV ilk model eklendi

$wb = $excel.ActiveWorkbook
$wsHcv = $wb.Worksheets.Item("HCV")

# HCV sheet: relabel the average-net-price column header from
# "ortalama_net_fiyat" to the new shared "net_fiyat" label.
$wsHcv.Range("B1").Value = "net_fiyat"

# Leave LCV_OTV's last selection (B5) alone; just make HCV the
# active/selected sheet with B1 selected, matching the saved view state.
$wsHcv.Activate()
$wsHcv.Range("B1").Select() | Out-Null
